$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Total Tcexecuted" (C) and "Passtest case" (D) columns entirely;
# remaining columns (Fail cases, Comments, Result After analysing) shift left.
$ws.Range("C:D").Delete()

# Update the active selection left after the edit.
$ws.Range("A12").Select()

$wb.Windows.Item(1).Width = 11130
